{"js": "// Replace the header date and the 25 division problems in the table\n// with the values from the target revision. Every source string in\n// this document is unique, so a simple search + replace per pair is\n// unambiguous and safe to apply in any order.\nconst replacements = [\n  [\"2026-02-02 Monday\", \"2026-02-03 Tuesday\"],\n  [\"582\u00f77=\", \"815\u00f75=\"],\n  [\"959\u00f77=\", \"747\u00f76=\"],\n  [\"857\u00f79=\", \"559\u00f78=\"],\n  [\"746\u00f75=\", \"149\u00f75=\"],\n  [\"313\u00f72=\", \"138\u00f75=\"],\n  [\"920\u00f75=\", \"437\u00f75=\"],\n  [\"668\u00f75=\", \"688\u00f75=\"],\n  [\"256\u00f73=\", \"562\u00f77=\"],\n  [\"871\u00f79=\", \"767\u00f78=\"],\n  [\"786\u00f79=\", \"665\u00f75=\"],\n  [\"894\u00f75=\", \"454\u00f77=\"],\n  [\"643\u00f72=\", \"315\u00f74=\"],\n  [\"471\u00f74=\", \"336\u00f72=\"],\n  [\"758\u00f78=\", \"642\u00f73=\"],\n  [\"763\u00f77=\", \"885\u00f79=\"],\n  [\"568\u00f78=\", \"361\u00f79=\"],\n  [\"941\u00f75=\", \"103\u00f75=\"],\n  [\"537\u00f78=\", \"936\u00f72=\"],\n  [\"517\u00f76=\", \"333\u00f74=\"],\n  [\"419\u00f77=\", \"226\u00f79=\"],\n  [\"698\u00f72=\", \"635\u00f75=\"],\n  [\"102\u00f79=\", \"808\u00f73=\"],\n  [\"701\u00f74=\", \"185\u00f74=\"],\n  [\"547\u00f74=\", \"532\u00f79=\"],\n  [\"418\u00f75=\", \"493\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the header date and the 25 division problems in the table\n# with the values from the target revision. Every source string in\n# this document is unique, so Find/Replace per pair is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @{old = \"2026-02-02 Monday\"; new = \"2026-02-03 Tuesday\"},\n  @{old = \"582\u00f77=\";            new = \"815\u00f75=\"},\n  @{old = \"959\u00f77=\";            new = \"747\u00f76=\"},\n  @{old = \"857\u00f79=\";            new = \"559\u00f78=\"},\n  @{old = \"746\u00f75=\";            new = \"149\u00f75=\"},\n  @{old = \"313\u00f72=\";            new = \"138\u00f75=\"},\n  @{old = \"920\u00f75=\";            new = \"437\u00f75=\"},\n  @{old = \"668\u00f75=\";            new = \"688\u00f75=\"},\n  @{old = \"256\u00f73=\";            new = \"562\u00f77=\"},\n  @{old = \"871\u00f79=\";            new = \"767\u00f78=\"},\n  @{old = \"786\u00f79=\";            new = \"665\u00f75=\"},\n  @{old = \"894\u00f75=\";            new = \"454\u00f77=\"},\n  @{old = \"643\u00f72=\";            new = \"315\u00f74=\"},\n  @{old = \"471\u00f74=\";            new = \"336\u00f72=\"},\n  @{old = \"758\u00f78=\";            new = \"642\u00f73=\"},\n  @{old = \"763\u00f77=\";            new = \"885\u00f79=\"},\n  @{old = \"568\u00f78=\";            new = \"361\u00f79=\"},\n  @{old = \"941\u00f75=\";            new = \"103\u00f75=\"},\n  @{old = \"537\u00f78=\";            new = \"936\u00f72=\"},\n  @{old = \"517\u00f76=\";            new = \"333\u00f74=\"},\n  @{old = \"419\u00f77=\";            new = \"226\u00f79=\"},\n  @{old = \"698\u00f72=\";            new = \"635\u00f75=\"},\n  @{old = \"102\u00f79=\";            new = \"808\u00f73=\"},\n  @{old = \"701\u00f74=\";            new = \"185\u00f74=\"},\n  @{old = \"547\u00f74=\";            new = \"532\u00f79=\"},\n  @{old = \"418\u00f75=\";            new = \"493\u00f75=\"}\n)\n\nforeach ($pair in $pairs) {\n  $rng = $d.Content\n  $null = $rng.Find.Execute($pair.old, $false, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)\n}\n"}
